# Fruta / hortaliza, semanal
# New weekly price records for "Pepino dulce" (Mapocho Venta Directa de
# Santiago) are inserted ahead of the existing rows 199-201, pushing all
# the subsequent rows down by three; the three rows that fall off the
# bottom of the table are appended again at the end (the table keeps a
# rolling window of the same length plus the new week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three fresh rows right before the current row 199, shifting
# rows 199:246 down to 202:249 and carrying the row-248/249 append that
# results from the sheet already ending right after row 246.
$ws.Range("A199:A201").EntireRow.Insert()

$newRows = @(
    @{ Row = 199; Fecha = 44798; Calidad = "Especial"; Volumen = 200; PMin = 16000; PMax = 16000; PProm = 16000; Precio = 889 },
    @{ Row = 200; Fecha = 44798; Calidad = "Primera";  Volumen = 250; PMin = 14000; PMax = 14000; PProm = 14000; Precio = 778 },
    @{ Row = 201; Fecha = 44798; Calidad = "Segunda";  Volumen = 220; PMin = 12000; PMax = 12000; PProm = 12000; Precio = 667 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 12
    $ws.Cells.Item($row, 2).Value = "Mapocho Venta Directa de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
